$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38: Just Give Him a Serum | Hi-Potion of Strength
$ws.Range("H38").Value = 988.7857
$ws.Range("I38").Value = 48.6
$ws.Range("J38").Value = 1511.1111
$ws.Range("K38").Value = 145.8
$ws.Range("L38").Value = 4533.3333
$ws.Range("M38").Value = 226.2
$ws.Range("N38").Value = -5277.3333

# Row 43: Growing Is Knowing | Growth Formula Gamma
$ws.Range("H43").Value = 4973.5557
$ws.Range("I43").Value = 3565.5
$ws.Range("K43").Value = 3565.5
$ws.Range("M43").Value = -3496.5

# Row 86: Filling in the Blanks | Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 2224.2222
$ws.Range("J86").Value = 1664
$ws.Range("L86").Value = 1664
$ws.Range("N86").Value = -3910

# Row 89: Ink into Antiquity (L) | Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 2224.2222
$ws.Range("J89").Value = 1664
$ws.Range("L89").Value = 8320
$ws.Range("N89").Value = -19552

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 1676.9
$ws.Range("I98").Value = 1676.9
$ws.Range("K98").Value = 1676.9
$ws.Range("M98").Value = -178.9000000000001

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 1676.9
$ws.Range("I122").Value = 1676.9
$ws.Range("K122").Value = 5030.700000000001
$ws.Range("M122").Value = -2580.700000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 2401.2632
$ws.Range("I32").Value = 2632.9375
$ws.Range("J32").Value = 1165.6666
$ws.Range("K32").Value = 2632.9375
$ws.Range("L32").Value = 1165.6666
$ws.Range("M32").Value = -2345.9375
$ws.Range("N32").Value = -1739.6666

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 3221.5557
$ws.Range("I61").Value = 3171
$ws.Range("J61").Value = 3398.5
$ws.Range("K61").Value = 3171
$ws.Range("L61").Value = 3398.5
$ws.Range("M61").Value = -2959
$ws.Range("N61").Value = -3822.5

# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Range("H102").Value = 2924.25
$ws.Range("I102").Value = 2879
$ws.Range("K102").Value = 2879
$ws.Range("M102").Value = -1257

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 3221.5557
$ws.Range("I136").Value = 3171
$ws.Range("J136").Value = 3398.5
$ws.Range("K136").Value = 9513
$ws.Range("L136").Value = 10195.5
$ws.Range("M136").Value = -6963
$ws.Range("N136").Value = -15295.5

$ws = $wb.Worksheets.Item("BSM")
# Row 64: With Bearings Straight | Mythrite Nugget
$ws.Range("H64").Value = 701
$ws.Range("J64").Value = 701
$ws.Range("L64").Value = 701
$ws.Range("N64").Value = -1151

# Row 67: Bearing the Brunt (L) | Mythrite Nugget
$ws.Range("H67").Value = 701
$ws.Range("J67").Value = 701
$ws.Range("L67").Value = 701
$ws.Range("N67").Value = -2261

# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 1475
$ws.Range("I86").Value = 1100
$ws.Range("J86").Value = 1600
$ws.Range("K86").Value = 1100
$ws.Range("L86").Value = 1600
$ws.Range("M86").Value = 23
$ws.Range("N86").Value = -3846

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 1475
$ws.Range("I89").Value = 1100
$ws.Range("J89").Value = 1600
$ws.Range("K89").Value = 5500
$ws.Range("L89").Value = 8000
$ws.Range("M89").Value = 116
$ws.Range("N89").Value = -19232

# Row 94: High Steal | High Steel Nugget
$ws.Range("H94").Value = 9966
$ws.Range("I94").Value = 9999
$ws.Range("K94").Value = 9999
$ws.Range("M94").Value = -9548

# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Range("H99").Value = 2619.1
$ws.Range("I99").Value = 1680.125
$ws.Range("K99").Value = 1680.125
$ws.Range("M99").Value = -182.125

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall | Elm Lumber
$ws.Range("H22").Value = 3364.3333
$ws.Range("I22").Value = 3364.3333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3364.3333
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -3014.3333
$ws.Range("N22").ClearContents()

# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("H105").Value = 1287.6666
$ws.Range("I105").Value = 1263
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 1263
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = 484
$ws.Range("N105").Value = -4794

$ws = $wb.Worksheets.Item("CUL")
# Row 58: Bread in the Clouds | La Noscean Toast
$ws.Range("H58").Value = 1350
$ws.Range("I58").Value = 1350
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 4050
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -3922
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 6500
$ws.Range("I126").Value = 6500
$ws.Range("K126").Value = 19500
$ws.Range("M126").Value = -17030

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 4509.1665
$ws.Range("I132").Value = 3986.5
$ws.Range("J132").Value = 5554.5
$ws.Range("K132").Value = 11959.5
$ws.Range("L132").Value = 16663.5
$ws.Range("M132").Value = -9429.5
$ws.Range("N132").Value = -21723.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 2099.25
$ws.Range("I7").Value = 2213.7144
$ws.Range("K7").Value = 2213.7144
$ws.Range("M7").Value = -2101.7144

# Row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws.Range("H55").Value = 2309.5557
$ws.Range("I55").Value = 2458
$ws.Range("J55").Value = 2124
$ws.Range("K55").Value = 2458
$ws.Range("L55").Value = 2124
$ws.Range("M55").Value = -2285
$ws.Range("N55").Value = -2470

# Row 68: You Could Say It's a Moving Target | Wyvern Leather
$ws.Range("H68").Value = 2615.625
$ws.Range("I68").Value = 2560.7144
$ws.Range("K68").Value = 2560.7144
$ws.Range("M68").Value = -1811.7144

# Row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Range("H71").Value = 2615.625
$ws.Range("I71").Value = 2560.7144
$ws.Range("K71").Value = 12803.572
$ws.Range("M71").Value = -9059.572

# Row 82: Trainin' the Neck | Dragon Leather
$ws.Range("H82").Value = 1313.75
$ws.Range("I82").Value = 1292.5714
$ws.Range("K82").Value = 1292.5714
$ws.Range("M82").Value = -931.5714

# Row 85: Training Is Only Skintight (L) | Dragon Leather
$ws.Range("H85").Value = 1313.75
$ws.Range("I85").Value = 1292.5714
$ws.Range("K85").Value = 1292.5714
$ws.Range("M85").Value = -44.57140000000004

# Row 100: Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 10143.556
$ws.Range("J100").Value = 18500
$ws.Range("L100").Value = 18500
$ws.Range("N100").Value = -19582

# Row 109: Band Substances | Smilodonskin Wristband
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 2099.25
$ws.Range("I126").Value = 2213.7144
$ws.Range("K126").Value = 6641.1432
$ws.Range("M126").Value = -4171.1432

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 4599.6665
$ws.Range("I132").Value = 4599.6665
$ws.Range("K132").Value = 13798.9995
$ws.Range("M132").Value = -11268.9995

# Row 134: Freezing Fingers | Crocodileskin Fingerless Gloves of Striking
$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 60000
$ws.Range("N134").Value = -70140

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 2360.3333
$ws.Range("I136").Value = 1508.6
$ws.Range("K136").Value = 4525.799999999999
$ws.Range("M136").Value = -1975.799999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display | Ruby Cotton Cloth
$ws.Range("H96").Value = 14
$ws.Range("I96").Value = 14
$ws.Range("K96").Value = 14
$ws.Range("M96").Value = 1359

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 6124.25
$ws.Range("I122").Value = 6124.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18372.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15922.75
$ws.Range("N122").ClearContents()

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 685.3333
$ws.Range("J126").Value = 999.5
$ws.Range("L126").Value = 2998.5
$ws.Range("N126").Value = -7938.5

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 1445.1
$ws.Range("I136").Value = 1445.1
$ws.Range("K136").Value = 4335.299999999999
$ws.Range("M136").Value = -1785.299999999999
